$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1907308377896613
$ws.Range("C2").Value = 0.5472370766488414
$ws.Range("J2").Value = 0.0106951871657754
$ws.Range("P2").Value = 0.1443850267379679
$ws.Range("S2").Value = 0.106951871657754
$ws.Range("B3").Value = 0.009584664536741214
$ws.Range("C3").Value = 0.03194888178913738
$ws.Range("J3").Value = 0.03194888178913738
$ws.Range("P3").Value = 0.7220447284345048
$ws.Range("S3").Value = 0.2044728434504792
$ws.Range("J4").Value = 0.0410958904109589
$ws.Range("O4").Value = 0.0136986301369863
$ws.Range("P4").Value = 0.6575342465753424
$ws.Range("S4").Value = 0.2876712328767123
$ws.Range("B6").Value = 0.06060606060606061
$ws.Range("D6").Value = 0.008080808080808081
$ws.Range("E6").Value = 0.00202020202020202
$ws.Range("F6").Value = 0.0505050505050505
$ws.Range("J6").Value = 0.2545454545454545
$ws.Range("O6").Value = 0.01212121212121212
$ws.Range("Q6").Value = 0.1656565656565657
$ws.Range("R6").Value = 0.08686868686868687
$ws.Range("S6").Value = 0.3595959595959596
$ws.Range("B7").Value = 0.1256281407035176
$ws.Range("D7").Value = 0.01507537688442211
$ws.Range("F7").Value = 0.07788944723618091
$ws.Range("J7").Value = 0.1055276381909548
$ws.Range("O7").Value = 0.01758793969849246
$ws.Range("Q7").Value = 0.1984924623115578
$ws.Range("R7").Value = 0.05025125628140704
$ws.Range("S7").Value = 0.4095477386934673
$ws.Range("B8").Value = 0.0824524312896406
$ws.Range("D8").Value = 0.01585623678646934
$ws.Range("E8").Value = 0.002114164904862579
$ws.Range("F8").Value = 0.06448202959830866
$ws.Range("J8").Value = 0.1025369978858351
$ws.Range("O8").Value = 0.02642706131078224
$ws.Range("Q8").Value = 0.200845665961945
$ws.Range("R8").Value = 0.1004228329809725
$ws.Range("S8").Value = 0.4048625792811839
$ws.Range("B9").Value = 0.09269662921348315
$ws.Range("D9").Value = 0.03089887640449438
$ws.Range("E9").Value = 0.002808988764044944
$ws.Range("F9").Value = 0.08146067415730338
$ws.Range("J9").Value = 0.1095505617977528
$ws.Range("O9").Value = 0.01966292134831461
$ws.Range("Q9").Value = 0.1882022471910112
$ws.Range("R9").Value = 0.09269662921348315
$ws.Range("S9").Value = 0.3820224719101123
$ws.Range("B10").Value = 0.1073049938093273
$ws.Range("D10").Value = 0.01650846058605035
$ws.Range("E10").Value = 0.0008254230293025176
$ws.Range("F10").Value = 0.0821295914156005
$ws.Range("J10").Value = 0.1064795707800248
$ws.Range("O10").Value = 0.02393726784977301
$ws.Range("Q10").Value = 0.2162608336772596
$ws.Range("R10").Value = 0.08419314898885678
$ws.Range("S10").Value = 0.3623607098638052
$ws.Range("G11").Value = 0.125
$ws.Range("J11").Value = 0.08854166666666667
$ws.Range("K11").Value = 0.1666666666666667
$ws.Range("L11").Value = 0.6024305555555556
$ws.Range("S11").Value = 0.01736111111111111
$ws.Range("G12").Value = 0.7320441988950276
$ws.Range("J12").Value = 0.1906077348066298
$ws.Range("K12").Value = 0.005524861878453038
$ws.Range("L12").Value = 0.03867403314917127
$ws.Range("S12").Value = 0.03314917127071823
$ws.Range("G13").Value = 0.6701030927835051
$ws.Range("J13").Value = 0.2989690721649484
$ws.Range("S13").Value = 0.03092783505154639
$ws.Range("F15").Value = 0.01867219917012448
$ws.Range("H15").Value = 0.1618257261410788
$ws.Range("I15").Value = 0.04356846473029045
$ws.Range("J15").Value = 0.3817427385892116
$ws.Range("K15").Value = 0.06639004149377593
$ws.Range("M15").Value = 0.008298755186721992
$ws.Range("N15").Value = 0.002074688796680498
$ws.Range("O15").Value = 0.08091286307053942
$ws.Range("S15").Value = 0.2365145228215768
$ws.Range("F16").Value = 0.02046783625730994
$ws.Range("H16").Value = 0.1754385964912281
$ws.Range("I16").Value = 0.08187134502923976
$ws.Range("J16").Value = 0.3976608187134503
$ws.Range("K16").Value = 0.1257309941520468
$ws.Range("M16").Value = 0.01754385964912281
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1286549707602339
$ws.Range("F17").Value = 0.01918976545842218
$ws.Range("H17").Value = 0.1748400852878465
$ws.Range("I17").Value = 0.08742004264392324
$ws.Range("J17").Value = 0.4339019189765458
$ws.Range("K17").Value = 0.09381663113006397
$ws.Range("M17").Value = 0.02665245202558635
$ws.Range("N17").Value = 0.001066098081023454
$ws.Range("O17").Value = 0.07995735607675906
$ws.Range("S17").Value = 0.08315565031982942
$ws.Range("F18").Value = 0.02284263959390863
$ws.Range("H18").Value = 0.2157360406091371
$ws.Range("I18").Value = 0.09390862944162437
$ws.Range("J18").Value = 0.3527918781725888
$ws.Range("K18").Value = 0.08629441624365482
$ws.Range("M18").Value = 0.02030456852791878
$ws.Range("N18").Value = 0.002538071065989848
$ws.Range("O18").Value = 0.07360406091370558
$ws.Range("S18").Value = 0.1319796954314721
$ws.Range("F19").Value = 0.01670146137787056
$ws.Range("H19").Value = 0.237160751565762
$ws.Range("I19").Value = 0.08225469728601252
$ws.Range("J19").Value = 0.3569937369519833
$ws.Range("K19").Value = 0.1144050104384134
$ws.Range("M19").Value = 0.02338204592901879
$ws.Range("N19").Value = 0.001670146137787056
$ws.Range("O19").Value = 0.07181628392484342
$ws.Range("S19").Value = 0.09561586638830898
